# This script applies a row-permutation update to worksheet "Artfynd".
# Rows 2-8 (species observation records) had their contents redistributed: each
# target row ends up holding the values that (in the prior version of the file)
# belonged to a different row, per the mapping: 2<-5, 3<-7, 4<-2, 5<-6, 6<-8, 7<-4, 8<-3.
# Only the cells whose value actually changes are touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2  (now holds the data formerly in row 5)
$ws.Range("A2").Value = 111863001
$ws.Range("B2").Value = 90332
$ws.Range("E2").Value = 4769
$ws.Range("F2").Value = "Svavelriska"
$ws.Range("G2").Value = "Lactarius scrobiculatus"
$ws.Range("H2").Value = "(Scop.:Fr.) Fr."
$ws.Range("P2").Value = "Charlottenberg, Upl"
$ws.Range("Q2").Value = 655217.6931657954
$ws.Range("R2").Value = 6634939.780080916
$ws.Range("Z2").Value = "10:47"
$ws.Range("AB2").Value = "10:47"

# Row 3  (now holds the data formerly in row 7)
$ws.Range("A3").Value = 111863218
$ws.Range("B3").Value = 90021
$ws.Range("E3").Value = 6031
$ws.Range("F3").Value = "Blomkålssvamp"
$ws.Range("G3").Value = "Sparassis crispa"
$ws.Range("H3").Value = "(Wulfen:Fr.) Fr."
$ws.Range("I3").Value = "1"
$ws.Range("J3").Value = "fruktkroppar"
$ws.Range("Q3").Value = 655137.9235184891
$ws.Range("R3").Value = 6634821.151011234
$ws.Range("Z3").Value = "10:53"
$ws.Range("AB3").Value = "10:53"
$ws.Range("AC3").ClearContents()

# Row 4  (now holds the data formerly in row 2)
$ws.Range("A4").Value = 111863402
$ws.Range("B4").Value = 90687
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 5964
$ws.Range("F4").Value = "Fjällig taggsvamp s.str."
$ws.Range("G4").Value = "Sarcodon imbricatus s.str."
$ws.Range("H4").Value = "(L.:Fr.) P.Karst."
$ws.Range("I4").Value = "1"
$ws.Range("Q4").Value = 655199.5794486763
$ws.Range("R4").Value = 6634769.85474884
$ws.Range("AC4").ClearContents()

# Row 5  (now holds the data formerly in row 6)
$ws.Range("A5").Value = 111863073
$ws.Range("B5").Value = 88899
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 3286
$ws.Range("F5").Value = "Flattoppad klubbsvamp"
$ws.Range("G5").Value = "Clavariadelphus truncatus"
$ws.Range("H5").Value = "(Quél.) Donk"
$ws.Range("I5").Value = "2"
$ws.Range("P5").Value = "Charlottenberg (Charlottenberg), Upl"
$ws.Range("Q5").Value = 655228.290648401
$ws.Range("R5").Value = 6634879.303300899
$ws.Range("Z5").Value = "10:50"
$ws.Range("AB5").Value = "10:50"

# Row 6  (now holds the data formerly in row 8)
$ws.Range("A6").Value = 111863288
$ws.Range("B6").Value = 85062
$ws.Range("E6").Value = 249278
$ws.Range("F6").Value = "Barrviolspindling"
$ws.Range("G6").Value = "Cortinarius harcynicus"
$ws.Range("H6").Value = "(Pers.) M.M.Moser"
$ws.Range("I6").Value = "1"
$ws.Range("Q6").Value = 655134.5683182024
$ws.Range("R6").Value = 6634792.815828164
$ws.Range("Z6").Value = "11:02"
$ws.Range("AB6").Value = "11:02"

# Row 7  (now holds the data formerly in row 4)
$ws.Range("A7").Value = 111863269
$ws.Range("B7").Value = 85062
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 249278
$ws.Range("F7").Value = "Barrviolspindling"
$ws.Range("G7").Value = "Cortinarius harcynicus"
$ws.Range("H7").Value = "(Pers.) M.M.Moser"
$ws.Range("I7").Value = "4"
$ws.Range("Q7").Value = 655135.2812587479
$ws.Range("R7").Value = 6634799.89438487
$ws.Range("Z7").Value = "11:02"
$ws.Range("AB7").Value = "11:02"
$ws.Range("AC7").Value = "4 ex i gräsglänta under gran och tall."

# Row 8  (now holds the data formerly in row 3)
$ws.Range("A8").Value = 111863040
$ws.Range("B8").Value = 90687
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 5964
$ws.Range("F8").Value = "Fjällig taggsvamp s.str."
$ws.Range("G8").Value = "Sarcodon imbricatus s.str."
$ws.Range("H8").Value = "(L.:Fr.) P.Karst."
$ws.Range("I8").Value = ""
$ws.Range("Q8").Value = 655235.4020021557
$ws.Range("R8").Value = 6634878.090185729
$ws.Range("Z8").Value = "10:49"
$ws.Range("AB8").Value = "10:49"
$ws.Range("AC8").Value = "Halv häxring, 3 m i diameter"
$ws.Range("J8").ClearContents()
